$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 25: new task 24
$ws.Range("A25").Value = "24. Add feature to make chance to change data in the time of review"
$ws.Range("B25").Value = "Hight"
$ws.Range("C25").Value = "Open"

# Update task 19's description to mention the colour change addition
$ws.Range("A20").Value = "19. Make feature to add word to ban list by clicking ( + change colour of words)"

# Fill in row 26: new task 25
$ws.Range("A26").Value = "25. Update algorithm to find data by using {0}{1}..{n} markers"
$ws.Range("B26").Value = "Hight"
$ws.Range("C26").Value = "Open"

# Update the active selection to reflect where the editor left off
$ws.Range("C19").Select()
